# Update the "Status" values (column E) for a subset of rows in Sheet1
# to reflect the latest counts, per commit:
# "fix: chroma compatibility with sqlite and system prompt for offersgen"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E12").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("E22").Value = 7
$ws.Range("E23").Value = 6
$ws.Range("E24").Value = 8
$ws.Range("E25").Value = 7
$ws.Range("E26").Value = 9
